$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.941.87"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.637.99"
$ws.Range("E4").Value = "  +0.34%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.70%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0639"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.60"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.82%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0795"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.865.17"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "1.625.93"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("E16").Value = "  -0.22%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "62.98"
$c.Style = "Normal"
$ws.Range("D18").Value = "25.986.51"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("E19").Value = "  +0.30%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "193.02"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.71%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.37"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.93"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  +0.96%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "144.25"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.129"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("E28").Value = "  -1.90%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.56"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.25"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0504"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.31"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("E35").Value = "  +1.54%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.901"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "1.136.54"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +0.34%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.48"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.91%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "99.39"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.798"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "1.774.73"
$ws.Range("E46").Value = "  +2.44%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "56.77"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("E49").Value = "  +0.28%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "7.69"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  -0.69%  "
